# ---------------------------------------------------------------------------
# InputBusOperadorCompleto.xlsx tutorial-update edit
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# 1. Rename sheets -----------------------------------------------------------
$wb.Worksheets.Item("Prefixes").Name             = "Prefix"
$wb.Worksheets.Item("PredicateObjectMaps").Name  = "Predicate_Object"
$wb.Worksheets.Item("Functions").Name            = "Function"

# 2. Reorder sheets: Subject now comes before Source -------------------------
$subject = $wb.Worksheets.Item("Subject")
$source  = $wb.Worksheets.Item("Source")
[void]$subject.Move($source)

# 3. Update the "Subject" sheet's C2:C5 URIs ----------------------------------
$wsSubject = $wb.Worksheets.Item("Subject")
$wsSubject.Range("C2").Value = "http://vocab.ciudadesabiertas.es/recurso/transporte/autobus/linea/{line_id}"
$wsSubject.Range("C3").Value = "http://vocab.ciudadesabiertas.es/recurso/transporte/autobus/presentacion/pres-{line_id}"
$wsSubject.Range("C4").Value = "http://vocab.ciudadesabiertas.es/recurso/transporte/autobus/operador/{agency_id}"
$wsSubject.Range("C5").Value = "http://vocab.ciudadesabiertas.es/recurso/transporte/autobus/ruta/{route_id}"
[void]$wsSubject.Activate()
[void]$wsSubject.Range("B12").Select()

# 4. Add a new predicate-object row entry (G13 = idRoute) on Predicate_Object
$wsPredObj = $wb.Worksheets.Item("Predicate_Object")
$wsPredObj.Range("G13").Value = "idRoute"
[void]$wsPredObj.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
[void]$wsPredObj.Range("E19").Select()

# 5. Update "Function" sheet headers (Predicate->Feature, Object->Value) -----
$wsFunction = $wb.Worksheets.Item("Function")
$wsFunction.Range("B1").Value = "Feature"
$wsFunction.Range("C1").Value = "Value"

# 6. Add hyperlink on Prefix!B8 (the vocab base URL) --------------------------
$wsPrefix = $wb.Worksheets.Item("Prefix")
$cellB8 = $wsPrefix.Range("B8")
$addr = $cellB8.Value2
[void]$wsPrefix.Hyperlinks.Add($cellB8, $addr)
[void]$wsPrefix.Activate()
[void]$wsPrefix.Range("B8").Select()

# 7. Final active sheet / window state: Function sheet, D21 selected ---------
[void]$wsFunction.Activate()
[void]$wsFunction.Range("D21").Select()

$win = $excel.ActiveWindow
$win.Left = 9520
$win.Top = 460
